# Insert a new data row at row 226 (shifting the existing rows 226-301 down
# to 227-302) and populate it with the new "Femacal de La Calera" Poroto
# verde observation dated 2022-01-27 (serial 44588) for Provincia de
# Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("226:226").Insert()

$ws.Cells.Item(226, 1).Value = 3
$ws.Cells.Item(226, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(226, 3).Value = "Coquimbo"
$ws.Cells.Item(226, 4).Value = 44588
$ws.Cells.Item(226, 5).Value = 5
$ws.Cells.Item(226, 6).Value = 100112031
$ws.Cells.Item(226, 7).Value = "Poroto verde"
$ws.Cells.Item(226, 8).Value = "Magnum"
$ws.Cells.Item(226, 9).Value = "Primera"
$ws.Cells.Item(226, 10).Value = 73
$ws.Cells.Item(226, 11).Value = 39000
$ws.Cells.Item(226, 12).Value = 40000
$ws.Cells.Item(226, 13).Value = 39479
$ws.Cells.Item(226, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(226, 15).Value = "Provincia de Santiago"
$ws.Cells.Item(226, 16).Value = 1579
$ws.Cells.Item(226, 17).Value = 25
$ws.Cells.Item(226, 18).Value = "Hortaliza"
